# Daily attendance processing - reverse the order of names in the
# "Recorded By" (column G) list for every session row on the active sheet.
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#          "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System"
#
# Cells whose value only contains a single name (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -le 1) { continue }

    $reversed = ""
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        if ($reversed -ne "") { $reversed += ", " }
        $reversed += $parts[$i]
    }

    if ($reversed -ne $text) {
        $cell.Value2 = $reversed
    }
}
